$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> @(Price (D) new value or $null if unchanged, Volume(1h) (E) new value)
$updates = @{
    2  = @("27.722.35", "  -0.88%  ")
    3  = @("1.758.30",  "  -0.50%  ")
    4  = @($null,       "  -0.18%  ")
    5  = @("324.34",    "  +0.80%  ")
    6  = @("1.000",     "  -0.14%  ")
    7  = @("0.4631",    "  +8.49%  ")
    8  = @("0.3611",    "  -0.27%  ")
    9  = @("0.07520",   "  +0.80%  ")
    10 = @("42.18",     "  -2.37%  ")
    11 = @("1.101",     "  +0.62%  ")
    12 = @($null,       "  -0.14%  ")
    13 = @("20.83",     "  -1.11%  ")
    14 = @("6.027",     "  -1.11%  ")
    15 = @("7.121",     "  -2.71%  ")
    16 = @("1.760.91",  "  -1.99%  ")
    17 = @("92.42",     "  +0.19%  ")
    18 = @("0.00001068","  +0.67%  ")
    19 = @("0.06402",   "  -0.11%  ")
    20 = @("0.9998",    "  -0.04%  ")
    21 = @("16.83",     "  -1.82%  ")
    22 = @("5.818",     "  -2.51%  ")
    23 = @("27.770.10", "  -0.72%  ")
    24 = @("11.29",     "  -0.27%  ")
    25 = @("2.104",     "  +0.14%  ")
    26 = @("164.89",    "  +3.91%  ")
    27 = @("20.41",     "  +0.97%  ")
    28 = @("1.962.56",  "  -1.38%  ")
    29 = @("2.085",     "  -3.30%  ")
    30 = @("126.75",    "  +0.98%  ")
    31 = @("1.067",     "  -7.85%  ")
    32 = @("0.09279",   "  +4.43%  ")
    33 = @($null,       "  -2.14%  ")
    34 = @("5.545",     "  -1.18%  ")
    35 = @($null,       "  -4.05%  ")
    36 = @($null,       "  -0.28%  ")
    37 = @("0.2106",    "  +0.05%  ")
    38 = @("0.06039",   "  +0.41%  ")
    39 = @("0.6371",    "  +0.24%  ")
    40 = @($null,       "  -0.87%  ")
    41 = @($null,       "  +1.68%  ")
    42 = @("1.379",     "  -1.65%  ")
    43 = @("7.846",     "  +0.42%  ")
    44 = @("13.37",     "  -0.32%  ")
    45 = @("0.5925",    "  +0.00%  ")
    46 = @("3.714",     "  +0.43%  ")
    47 = @("123.60",    "  +0.88%  ")
    48 = @("1.956",     "  -2.46%  ")
    49 = @("1.151",     "  -2.78%  ")
    50 = @("0.06874",   "  +0.15%  ")
    51 = @("72.51",     "  -2.27%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    if ($priceVal -ne $null) {
        # Force the Price cell to remain plain text (it looks numeric, e.g.
        # "324.34" or "1.000", and Excel would otherwise coerce it into a
        # number and drop significant formatting such as trailing zeros).
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $priceVal
        $dCell.Style = "Normal"
    }

    # The Volume(1h) text (e.g. "  -0.88%  ") is never number-like, so a
    # plain assignment keeps it as text without touching cell style.
    $ws.Cells.Item($row, 5).Value = $volVal
}
